$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5", "D6", "D8", "D13", "D14", "D15", "D17", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D46", "D47", "D49", "D50", "D51")
foreach ($cellAddr in $textCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

$ws.Range("D2").Value = '57.367.84'
$ws.Range("E2").Value = '  -4.07%  '
$ws.Range("D3").Value = '2.909.92'
$ws.Range("E3").Value = '  -2.58%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '546.69'
$ws.Range("E5").Value = '  -4.06%  '
$ws.Range("D6").Value = '129.36'
$ws.Range("E6").Value = '  +3.07%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = '0.510'
$ws.Range("E8").Value = '  +2.03%  '
$ws.Range("D9").Value = '2.903.12'
$ws.Range("E9").Value = '  -2.69%  '
$ws.Range("E10").Value = '  -3.31%  '
$ws.Range("E11").Value = '  -6.71%  '
$ws.Range("E12").Value = '  +1.30%  '
$ws.Range("D13").Value = '0.0000219'
$ws.Range("E13").Value = '  +0.44%  '
$ws.Range("D14").Value = '32.69'
$ws.Range("E14").Value = '  +0.62%  '
$ws.Range("D15").Value = '0.120'
$ws.Range("E15").Value = '  +0.23%  '
$ws.Range("D16").Value = '3.394.72'
$ws.Range("E16").Value = '  -2.42%  '
$ws.Range("D17").Value = '6.81'
$ws.Range("E17").Value = '  +5.99%  '
$ws.Range("D18").Value = '2.910.82'
$ws.Range("E18").Value = '  -2.16%  '
$ws.Range("D19").Value = '57.401.27'
$ws.Range("E19").Value = '  -4.09%  '
$ws.Range("D20").Value = '416.04'
$ws.Range("E20").Value = '  -2.32%  '
$ws.Range("D21").Value = '13.07'
$ws.Range("E21").Value = '  -0.12%  '
$ws.Range("D22").Value = '0.680'
$ws.Range("E22").Value = '  +2.04%  '
$ws.Range("D23").Value = '6.92'
$ws.Range("E23").Value = '  -1.53%  '
$ws.Range("D24").Value = '13.00'
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").Value = '79.51'
$ws.Range("E25").Value = '  +0.52%  '
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("E28").Value = '  -2.69%  '
$ws.Range("D29").Value = '7.40'
$ws.Range("E29").Value = '  +1.93%  '
$ws.Range("E30").Value = '  +1.63%  '
$ws.Range("D31").Value = '25.07'
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("D32").Value = '5.89'
$ws.Range("E32").Value = '  -4.24%  '
$ws.Range("D33").Value = '0.0957'
$ws.Range("E33").Value = '  +2.36%  '
$ws.Range("D34").Value = '5.61'
$ws.Range("E34").Value = '  +0.17%  '
$ws.Range("D35").Value = '0.932'
$ws.Range("E35").Value = '  +0.42%  '
$ws.Range("D36").Value = '2.05'
$ws.Range("E36").Value = '  +0.31%  '
$ws.Range("D37").Value = '47.80'
$ws.Range("E37").Value = '  -4.26%  '
$ws.Range("D38").Value = '8.66'
$ws.Range("E38").Value = '  +4.13%  '
$ws.Range("D39").Value = '0.0₃0670'
$ws.Range("E39").Value = '  +2.16%  '
$ws.Range("D40").Value = '2.53'
$ws.Range("E40").Value = '  +2.76%  '
$ws.Range("D41").Value = '0.106'
$ws.Range("E41").Value = '  -1.50%  '
$ws.Range("D42").Value = '0.0343'
$ws.Range("E42").Value = '  -2.88%  '
$ws.Range("D43").Value = '371.40'
$ws.Range("E43").Value = '  -2.37%  '
$ws.Range("D44").Value = '2.658.43'
$ws.Range("E44").Value = '  -0.24%  '
$ws.Range("D46").Value = '121.64'
$ws.Range("E46").Value = '  +1.33%  '
$ws.Range("D47").Value = '0.236'
$ws.Range("E47").Value = '  +0.55%  '
$ws.Range("E48").Value = '  +1.63%  '
$ws.Range("D49").Value = '1.95'
$ws.Range("E49").Value = '  -1.97%  '
$ws.Range("D50").Value = '23.04'
$ws.Range("E50").Value = '  -2.54%  '
$ws.Range("D51").Value = '1.99'
$ws.Range("E51").Value = '  -0.03%  '
